{"js": "// \"fixed text smear, added skybox, added text scroll\"\n//\n// This script performs the text-smear fix described by the diff:\n//  1) \"I can start working on the things I actually want to now.\" was split\n//     across three runs (with <w:proofErr> gramStart/gramEnd markers sitting\n//     between them). Re-set the paragraph's text so it collapses back into a\n//     single clean run.\n//  2) Likewise for \"I could be doing something so much more productive with\n//     my time. Like reading a book, learning something actually useful. But\n//     I choose this\u2026\".\n//  3) Two new paragraphs (\"I really don't know what else to do with this\n//     game so I'm just gonna turn it off for you now.\" and \"Sorry.\") are\n//     added right after the (pre-existing) blank paragraph that follows that\n//     same \"I could be doing something\u2026\" paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the two paragraphs that need their split runs collapsed, and the\n// blank paragraph right after the second one (that's where the new lines\n// get inserted).\nlet wantNowIndex = -1;\nlet wantChooseIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (wantNowIndex === -1 && t.indexOf(\"I can start working on the things I actually want to\") === 0) {\n    wantNowIndex = i;\n  }\n  if (wantChooseIndex === -1 && t.indexOf(\"I could be doing something so much more productive\") === 0) {\n    wantChooseIndex = i;\n  }\n}\n\nif (wantNowIndex === -1 || wantChooseIndex === -1) {\n  throw new Error(\"Could not locate the target paragraphs to fix.\");\n}\n\n// 1) Re-write the text of each paragraph as a single contiguous string \u2014\n// Word.Paragraph.insertText(text, \"Replace\") replaces the paragraph's whole\n// range (and the w:proofErr-separated runs it held) with one plain run.\nitems[wantNowIndex].insertText(\n  \"I can start working on the things I actually want to now.\",\n  \"Replace\"\n);\n\nitems[wantChooseIndex].insertText(\n  \"I could be doing something so much more productive with my time. Like reading a book, learning something actually useful. But I choose this\\u2026\",\n  \"Replace\"\n);\n\n// 2) Insert the two new paragraphs after the blank paragraph that directly\n// follows the \"I could be doing something\u2026\" paragraph.\nconst blankAfterChoose = items[wantChooseIndex + 1];\nconst scrollLine1 = blankAfterChoose.insertParagraph(\n  \"I really don\\u2019t know what else to do with this game so I\\u2019m just gonna turn it off for you now.\",\n  \"After\"\n);\nscrollLine1.insertParagraph(\"Sorry.\", \"After\");\n\nawait context.sync();\n", "ps1": "# \"fixed text smear, added skybox, added text scroll\"\n#\n# 1) \"I can start working on the things I actually want to now.\" and\n#    \"I could be doing something so much more productive with my time. Like\n#    reading a book, learning something actually useful. But I choose this\u2026\"\n#    were each smeared across three runs with <w:proofErr> gramStart/gramEnd\n#    markers in between. Find/replace the whole sentence (Replace:=2,\n#    wdReplaceAll) so Word collapses it back into a single clean run.\n# 2) Insert two new paragraphs (\"I really don't know what else to do with\n#    this game so I'm just gonna turn it off for you now.\" and \"Sorry.\")\n#    right after the blank paragraph that follows the \"I could be doing\n#    something\u2026\" paragraph.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$idxNow = -1\n$idxChoose = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($idxNow -eq -1 -and $t.StartsWith(\"I can start working on the things I actually want to\")) {\n        $idxNow = $i\n    }\n    if ($idxChoose -eq -1 -and $t.StartsWith(\"I could be doing something so much more productive\")) {\n        $idxChoose = $i\n    }\n}\n\nif ($idxNow -eq -1 -or $idxChoose -eq -1) {\n    throw \"Could not locate the target paragraphs to fix.\"\n}\n\n# --- 1a) collapse the \"...actually want to now.\" run-smear ---------------\n$nowText = \"I can start working on the things I actually want to now.\"\n$rng = $d.Paragraphs.Item($idxNow).Range\n$find = $rng.Find\n$find.Text = $nowText\n$find.Replacement.Text = $nowText\n$null = $find.Execute($nowText, $false, $false, $false, $false, $false, $true, 1, $false, $nowText, 2)\n\n# --- 1b) collapse the \"...But I choose this\u2026\" run-smear -------------------\n$chooseText = \"I could be doing something so much more productive with my time. Like reading a book, learning something actually useful. But I choose this\" + [char]0x2026\n$rng2 = $d.Paragraphs.Item($idxChoose).Range\n$find2 = $rng2.Find\n$find2.Text = $chooseText\n$find2.Replacement.Text = $chooseText\n$null = $find2.Execute($chooseText, $false, $false, $false, $false, $false, $true, 1, $false, $chooseText, 2)\n\n# --- 2) add the two new \"text scroll\" paragraphs after the blank line -----\n$blank = $d.Paragraphs.Item($idxChoose + 1)\n$tail = $blank.Range\n$tail.Collapse(0)   # wdCollapseEnd\n$tail.InsertParagraphAfter()\n$tail.Collapse(0)\n$tail.InsertAfter(\"I really don\" + [char]0x2019 + \"t know what else to do with this game so I\" + [char]0x2019 + \"m just gonna turn it off for you now.\")\n\n$tail2 = $d.Paragraphs.Item($idxChoose + 2).Range\n$tail2.Collapse(0)\n$tail2.InsertParagraphAfter()\n$tail2.Collapse(0)\n$tail2.InsertAfter(\"Sorry.\")\n"}
